# Auto-generated edits applying scheduled runner updates to Ixion_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 225.66667
$ws.Range("I12").Value = 66.2
$ws.Range("J12").Value = 425
$ws.Range("K12").Value = 66.2
$ws.Range("L12").Value = 425
$ws.Range("M12").Value = 103.8
$ws.Range("N12").Value = -765

# Row 137
$ws.Range("H137").Value = 1517.3334
$ws.Range("I137").Value = 1206.0358
$ws.Range("J137").Value = 3260.6
$ws.Range("K137").Value = 3618.1074
$ws.Range("L137").Value = 9781.799999999999
$ws.Range("M137").Value = -1068.1074
$ws.Range("N137").Value = -14881.8

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1471.5454
$ws.Range("I61").Value = 1414.4445
$ws.Range("J61").Value = 1728.5
$ws.Range("K61").Value = 1414.4445
$ws.Range("L61").Value = 1728.5
$ws.Range("M61").Value = -1202.4445
$ws.Range("N61").Value = -2152.5

# Row 74
$ws.Range("H74").Value = 1345.0652
$ws.Range("I74").Value = 1186.1316
$ws.Range("K74").Value = 1186.1316
$ws.Range("M74").Value = -312.1315999999999

# Row 77
$ws.Range("H77").Value = 1345.0652
$ws.Range("I77").Value = 1186.1316
$ws.Range("K77").Value = 5930.657999999999
$ws.Range("M77").Value = -1562.657999999999

# Row 122
$ws.Range("H122").Value = 2138520.8
$ws.Range("I122").Value = 2332659
$ws.Range("K122").Value = 6997977
$ws.Range("M122").Value = -6995527

# Row 132
$ws.Range("H132").Value = 4192.0186
$ws.Range("I132").Value = 1257.8182
$ws.Range("J132").Value = 17102.5
$ws.Range("K132").Value = 3773.4546
$ws.Range("L132").Value = 51307.5
$ws.Range("M132").Value = -1243.4546
$ws.Range("N132").Value = -56367.5

# Row 136
$ws.Range("H136").Value = 1471.5454
$ws.Range("I136").Value = 1414.4445
$ws.Range("J136").Value = 1728.5
$ws.Range("K136").Value = 4243.333500000001
$ws.Range("L136").Value = 5185.5
$ws.Range("M136").Value = -1693.333500000001
$ws.Range("N136").Value = -10285.5

$ws = $wb.Worksheets.Item("BSM")
# Row 40
$ws.Range("H40").Value = 29000
$ws.Range("J40").Value = 29000
$ws.Range("L40").Value = 29000
$ws.Range("N40").Value = -29530

# Row 96
$ws.Range("H96").Value = 38000
$ws.Range("J96").Value = 38000
$ws.Range("L96").Value = 38000
$ws.Range("N96").Value = -43492

# Row 134
$ws.Range("H134").Value = 2423.1177
$ws.Range("I134").Value = 2181.0908
$ws.Range("J134").Value = 2866.8333
$ws.Range("K134").Value = 6543.2724
$ws.Range("L134").Value = 8600.499899999999
$ws.Range("M134").Value = -4008.2724
$ws.Range("N134").Value = -13670.4999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 10323.2
$ws.Range("I31").Value = 3322.4
$ws.Range("J31").Value = 13823.6
$ws.Range("K31").Value = 3322.4
$ws.Range("L31").Value = 13823.6
$ws.Range("M31").Value = -3027.4
$ws.Range("N31").Value = -14413.6

# Row 34
$ws.Range("H34").Value = 10323.2
$ws.Range("I34").Value = 3322.4
$ws.Range("J34").Value = 13823.6
$ws.Range("K34").Value = 3322.4
$ws.Range("L34").Value = 13823.6
$ws.Range("M34").Value = -3120.4
$ws.Range("N34").Value = -14227.6

# Row 58
$ws.Range("H58").Value = 999.2222
$ws.Range("I58").Value = 745.2857
$ws.Range("J58").Value = 1888
$ws.Range("K58").Value = 745.2857
$ws.Range("L58").Value = 1888
$ws.Range("M58").Value = -542.2857
$ws.Range("N58").Value = -2294

# Row 132
$ws.Range("H132").Value = 2255.64
$ws.Range("I132").Value = 2047.6666
$ws.Range("J132").Value = 2790.4285
$ws.Range("K132").Value = 6142.9998
$ws.Range("L132").Value = 8371.2855
$ws.Range("M132").Value = -3612.9998
$ws.Range("N132").Value = -13431.2855

# Row 134
$ws.Range("H134").Value = 1799.4166
$ws.Range("I134").Value = 1984.7142
$ws.Range("J134").Value = 1540
$ws.Range("K134").Value = 5954.142599999999
$ws.Range("L134").Value = 4620
$ws.Range("M134").Value = -3419.142599999999
$ws.Range("N134").Value = -9690

# Row 136
$ws.Range("H136").Value = 999.2222
$ws.Range("I136").Value = 745.2857
$ws.Range("J136").Value = 1888
$ws.Range("K136").Value = 2235.8571
$ws.Range("L136").Value = 5664
$ws.Range("M136").Value = 314.1428999999998
$ws.Range("N136").Value = -10764

$ws = $wb.Worksheets.Item("CUL")
# Row 117
$ws.Range("H117").Value = 15879666
$ws.Range("J117").Value = 17545918
$ws.Range("L117").Value = 52637754
$ws.Range("N117").Value = -52644638

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5540.592
$ws.Range("I70").Value = 5578.081
$ws.Range("J70").Value = 5425
$ws.Range("K70").Value = 5578.081
$ws.Range("L70").Value = 5425
$ws.Range("M70").Value = -5308.081
$ws.Range("N70").Value = -5965

# Row 73
$ws.Range("H73").Value = 5540.592
$ws.Range("I73").Value = 5578.081
$ws.Range("J73").Value = 5425
$ws.Range("K73").Value = 5578.081
$ws.Range("L73").Value = 5425
$ws.Range("M73").Value = -4642.081
$ws.Range("N73").Value = -7297

# Row 102
$ws.Range("H102").Value = 1188.0834
$ws.Range("I102").Value = 1113
$ws.Range("J102").Value = 2014
$ws.Range("K102").Value = 1113
$ws.Range("L102").Value = 2014
$ws.Range("M102").Value = 509
$ws.Range("N102").Value = -5258

# Row 126
$ws.Range("H126").Value = 5744.68
$ws.Range("I126").Value = 7539.4707
$ws.Range("J126").Value = 1930.75
$ws.Range("K126").Value = 22618.4121
$ws.Range("L126").Value = 5792.25
$ws.Range("M126").Value = -20148.4121
$ws.Range("N126").Value = -10732.25

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 35716620
$ws.Range("I40").Value = 41668772
$ws.Range("J40").Value = 3698.75
$ws.Range("K40").Value = 41668772
$ws.Range("L40").Value = 3698.75
$ws.Range("M40").Value = -41668636
$ws.Range("N40").Value = -3970.75

# Row 132
$ws.Range("H132").Value = 27086060
$ws.Range("I132").Value = 48149916
$ws.Range("K132").Value = 144449748
$ws.Range("M132").Value = -144447218

# Row 136
$ws.Range("H136").Value = 4061.375
$ws.Range("I136").Value = 2798.16
$ws.Range("J136").Value = 8572.857
$ws.Range("K136").Value = 8394.48
$ws.Range("L136").Value = 25718.571
$ws.Range("M136").Value = -5844.48
$ws.Range("N136").Value = -30818.571

$ws = $wb.Worksheets.Item("WVR")
# Row 99
$ws.Range("H99").Value = 21666.666
$ws.Range("I99").Value = 20000
$ws.Range("K99").Value = 20000
$ws.Range("M99").Value = -17005

# Row 132
$ws.Range("H132").Value = 1494.871
$ws.Range("I132").Value = 1010.2083
$ws.Range("J132").Value = 3156.5715
$ws.Range("K132").Value = 3030.6249
$ws.Range("L132").Value = 9469.7145
$ws.Range("M132").Value = -500.6248999999998
$ws.Range("N132").Value = -14529.7145

# Row 136
$ws.Range("H136").Value = 4367.5
$ws.Range("I136").Value = 4956.6665
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 14869.9995
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -12319.9995
$ws.Range("N136").Value = -12900

Write-Host "Applied all Ixion_Profits edits"
